$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.442.26"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "3.363.93"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'257.99"
$ws.Range("E5").Value = "  -0.05%  "
$ws.Range("D6").Value = "'665.63"
$ws.Range("E6").Value = "  +6.14%  "
$ws.Range("E7").Value = "  +9.95%  "
$ws.Range("D8").Value = "'0.471"
$ws.Range("E8").Value = "  +21.88%  "
$ws.Range("D9").Value = "'1.08"
$ws.Range("E9").Value = "  +25.18%  "
$ws.Range("E10").Value = "  -0.01%  "
$ws.Range("D11").Value = "3.362.66"
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("D12").Value = "'0.215"
$ws.Range("E12").Value = "  +8.32%  "
$ws.Range("D13").Value = "'42.35"
$ws.Range("E13").Value = "  +14.83%  "
$ws.Range("D14").Value = "'0.0000273"
$ws.Range("E14").Value = "  +10.52%  "
$ws.Range("D15").Value = "98.568.69"
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").Value = "3.991.59"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").Value = "'5.69"
$ws.Range("E17").Value = "  +3.78%  "
$ws.Range("D18").Value = "3.359.80"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "'7.68"
$ws.Range("E19").Value = "  +27.14%  "
$ws.Range("D20").Value = "'16.79"
$ws.Range("E20").Value = "  +10.65%  "
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("D22").Value = "'532.50"
$ws.Range("E22").Value = "  +9.75%  "
$ws.Range("D23").Value = "'10.56"
$ws.Range("E23").Value = "  +12.47%  "
$ws.Range("D24").Value = "'0.0000221"
$ws.Range("E24").Value = "  +5.53%  "
$ws.Range("D25").Value = "'0.438"
$ws.Range("E25").Value = "  +54.45%  "
$ws.Range("D26").Value = "'102.27"
$ws.Range("E26").Value = "  +15.34%  "
$ws.Range("D27").Value = "'6.19"
$ws.Range("E27").Value = "  +10.13%  "
$ws.Range("D28").Value = "'12.59"
$ws.Range("E28").Value = "  +6.21%  "
$ws.Range("D29").Value = "3.544.62"
$ws.Range("E29").Value = "  -0.12%  "
$ws.Range("D30").Value = "'0.147"
$ws.Range("E30").Value = "  +7.67%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").Value = "'11.03"
$ws.Range("E32").Value = "  +14.75%  "
$ws.Range("D33").Value = "'0.190"
$ws.Range("E33").Value = "  -1.51%  "
$ws.Range("E34").Value = "  -0.09%  "
$ws.Range("D35").Value = "'29.39"
$ws.Range("E35").Value = "  +5.50%  "
$ws.Range("D36").Value = "'0.540"
$ws.Range("E36").Value = "  +17.74%  "
$ws.Range("D37").Value = "'7.86"
$ws.Range("E37").Value = "  +8.14%  "
$ws.Range("E38").Value = "  +8.88%  "
$ws.Range("E39").Value = "  +5.45%  "
$ws.Range("D40").Value = "'527.48"
$ws.Range("E40").Value = "  +6.08%  "
$ws.Range("E41").Value = "  +6.32%  "
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("B43").Value = "MantraDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D43").Value = "'3.80"
$ws.Range("E43").Value = "  +2.67%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0431"
$ws.Range("E44").Value = "  +32.34%  "
$ws.Range("D45").Value = "'3.44"
$ws.Range("E45").Value = "  +4.70%  "
$ws.Range("D46").Value = "'0.828"
$ws.Range("E46").Value = "  +6.21%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "'2.07"
$ws.Range("E48").Value = "  +7.59%  "
$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").Value = "'5.17"
$ws.Range("E49").Value = "  +12.75%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D50").Value = "'7.91"
$ws.Range("E50").Value = "  +18.58%  "
$ws.Range("D51").Value = "'50.97"
$ws.Range("E51").Value = "  +11.81%  "
